$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set column H (Absent) to 1 for rows 3 through 18
$ws.Range("H3:H18").Value = 1
